# Add the new "RecordCreditNote" worksheet as the last tab (this also
# activates it, matching the workbook's new activeTab and the sheet's own
# tabSelected flag, and correspondingly clears the previous active sheet's
# selected-tab flag).
$wb = $excel.ActiveWorkbook
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "RecordCreditNote"

# Column widths (characters) roughly matching the source layout.
$ws.Columns.Item(1).ColumnWidth = 35.084
$ws.Columns.Item(2).ColumnWidth = 31.251
$ws.Columns.Item(3).ColumnWidth = 23.75
$ws.Columns.Item(4).ColumnWidth = 29.75
$ws.Columns.Item(5).ColumnWidth = 36.084
$ws.Columns.Item(6).ColumnWidth = 34.75

# Header row.
$ws.Cells.Item(1, 1).Value = "RECORDCREDITNOTE_UBIN_ITEMNUMBER"
$ws.Cells.Item(1, 2).Value = "RECORDCREDITNOTE_WEARER_NAME"
$ws.Cells.Item(1, 3).Value = "RECORDCREDITNOTE_NOTES"
$ws.Cells.Item(1, 4).Value = "RECORDCREDITNOTE_OBSERVATION"
$ws.Cells.Item(1, 5).Value = "RECORDCREDITNOTE_START_DATE_FORMAT"
$ws.Cells.Item(1, 6).Value = "RECORDCREDITNOTE_END_DATE_FORMAT"

# Data row.
$ws.Cells.Item(2, 1).Value = "Testing"
$ws.Cells.Item(2, 2).Value = "Wearer name from excel"
$ws.Cells.Item(2, 3).Value = "Credit notes from excel sheet"
$ws.Cells.Item(2, 4).Value = "Observation from excel"

$ws.Cells.Item(2, 5).Value = "16/10/2020"
$ws.Cells.Item(2, 5).NumberFormat = "DD/MM/YY"
$ws.Cells.Item(2, 6).Value = "17/10/2020"
$ws.Cells.Item(2, 6).NumberFormat = "DD/MM/YY"
